$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins, Losses, Ties in AD1:AF1,
# matching the style used by the other header cells (A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("AD1:AF1").VerticalAlignment = -4160    # xlTop
$ws.Range("AD1:AF1").Borders.LineStyle = 1
$ws.Range("AD1:AF1").Borders.Weight = 2

# Fill in the team record (Wins/Losses/Ties) for every data row (2-47).
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 82  # AD
    $ws.Cells.Item($r, 31).Value = 80  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
